$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 3.7107675
$ws.Range("H2").Value = 7.421535
$ws.Range("I2").Value = 0.4297212203365021
$ws.Range("J2").Value = 0.3540438930304464
$ws.Range("Q2").Value = 1.033995468495
$ws.Range("R2").Value = 6.203972810970001
$ws.Range("S2").Value = 0.4297212203365021
$ws.Range("T2").Value = 0.3540438930304464

$ws.Range("G3").Value = 0.8756340000000001
$ws.Range("I3").Value = 0.1014018019313074
$ws.Range("J3").Value = 0.1253162062416287
$ws.Range("Q3").Value = 0.2439930790760001
$ws.Range("S3").Value = 0.1014018019313074
$ws.Range("T3").Value = 0.1253162062416287

$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 1.048369666666667
$ws.Range("H4").Value = 3.145109
$ws.Range("I4").Value = 0.1214052598347301
$ws.Range("J4").Value = 0.1500372408625836
$ws.Range("Q4").Value = 0.2921254119642223
$ws.Range("R4").Value = 2.629128707678
$ws.Range("S4").Value = 0.1214052598347301
$ws.Range("T4").Value = 0.1500372408625836

$ws.Range("G5").Value = 1.2329145
$ws.Range("H5").Value = 2.465829
$ws.Range("I5").Value = 0.1427762648860562
$ws.Range("J5").Value = 0.1176322281990683
$ws.Range("Q5").Value = 0.3435483376530001
$ws.Range("R5").Value = 2.061290025918
$ws.Range("S5").Value = 0.1427762648860562
$ws.Range("T5").Value = 0.1176322281990683

$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.9623080000000001
$ws.Range("H6").Value = 2.886924
$ws.Range("I6").Value = 0.1114389861664948
$ws.Range("J6").Value = 0.1377205405408758
$ws.Range("Q6").Value = 0.2681445580453334
$ws.Range("R6").Value = 2.413301022408
$ws.Range("S6").Value = 0.1114389861664948
$ws.Range("T6").Value = 0.1377205405408758

$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.8052966666666667
$ws.Range("H7").Value = 2.41589
$ws.Range("I7").Value = 0.0932564668449094
$ws.Range("J7").Value = 0.1152498911253973
$ws.Range("Q7").Value = 0.2243937687088889
$ws.Range("R7").Value = 2.01954391838
$ws.Range("S7").Value = 0.0932564668449094
$ws.Range("T7").Value = 0.1152498911253973
